$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "SCRIPT/T01P02A/um1103.ssb"
$ws.Range("B6").Value = 152
$ws.Range("C6").Value = " That stink is finally fading…"
$ws.Range("D6").Value = " Эта вонь наконец-то\nрассеивается…"
$ws.Range("E6").Value = " Üóà âïîû îàëïîåø-óï\nñàòòåéâàåóòÿ…"

$ws.Rows.Item(6).RowHeight = 43.2

$ws.Range("D1").Select() | Out-Null
